$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new site ("bidaha") was added to the alphabetically-sorted bonus list.
# It sorts between "bewins" (row 115) and "bonisa" (row 116), so insert a
# fresh row at 116 and push everything from there on down by one.
$ws.Rows("116:116").Insert()

$ws.Range("A116").Value = "bidaha"
$ws.Range("B116").Value = "Maks 10k çekim"
$ws.Range("C116").Value = "yatırımsız"

# Leave the selection where the author ended up after making the edit.
$ws.Range("B15").Select()
